# "chef mapping DB design"
# Redesign the "Order to chef" sheet: replace the old ChefOrder /
# MapOrderToChef two-table layout (columns A + C) with a new
# ChefOrderMapping / ChefOrderCancelDetail two-table layout (columns A + B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order to chef")

# --- clear the old column C table (keep the header style on C1) ---
$ws.Range("C1:C9").ClearContents()

# --- column A: ChefOrderMapping field list ---
$ws.Range("A1").Value  = "ChefOrderMapping"
$ws.Range("A2").Value  = "MappingId"
$ws.Range("A3").Value  = "ChefId"
$ws.Range("A4").Value  = "OrderDetailId"
$ws.Range("A5").Value  = "Quantity"
$ws.Range("A6").Value  = "UserId"
$ws.Range("A7").Value  = "OrderGivenTime"
$ws.Range("A8").Value  = "ExpectedPickupTime"
$ws.Range("A9").Value  = "ActualPickupTime"
$ws.Range("A10").Value = "CreatedBy"
$ws.Range("A11").Value = "CreatedOn"
$ws.Range("A12").Value = "ModifiedBy"
$ws.Range("A13").Value = "ModifiedOn"
$ws.Range("A14").Value = "Status"

# --- column B: ChefOrderCancelDetail field list (new column) ---
$ws.Range("B1").Value  = "ChefOrderCancelDetail"
$ws.Range("B2").Value  = "OrderCanceld"
$ws.Range("B3").Value  = "ChefId"
$ws.Range("B4").Value  = "OrderDetailId"
$ws.Range("B5").Value  = "Quantity"
$ws.Range("B6").Value  = "Remarks"
$ws.Range("B7").Value  = "CreatedBy"
$ws.Range("B8").Value  = "CreatedOn"
$ws.Range("B9").Value  = "ModifiedBy"
$ws.Range("B10").Value = "ModifiedOn"
$ws.Range("B11").Value = "Status"

# header row (A1, B1) is bold, matching the other header cells already on the sheet
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true

# give the new column B a best-fit-like width
$ws.Columns.Item(2).ColumnWidth = 19.4

# the Customer sheet's selection moves too (reviewer was looking at the
# UserAddressDetails block while this edit was made) -- set it without
# leaving that sheet active, since "Order to chef" is the tab in focus
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("C10:C14").Select()

# "Order to chef" stays the active/selected tab; dimension now runs
# through row 14 / column C, so refresh the view selection there too
$ws.Activate()
$ws.Range("C14").Select()

# --- sharedStrings for the workbook now carry eight new strings used above
#     and drop the eight that were removed (ChefOrder / ChefOrderId /
#     MapOrderID / ChefDeliveredDateTime / MapOrderToChef /
#     AssignedPickUpDate / AssignedPickUpTime / OrderGivenDatetime) — this
#     happens automatically as a consequence of the cell edits above.
